$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.383.85'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.848.72'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6277'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07618'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07739'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.031'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6797'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001056'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.19'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.160'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '29.409.67'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9996'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.484'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '158.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.1386'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.405'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('E27').Value = '  +8.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.460'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05606'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.113'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.068'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.6999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.588'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01803'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '1.230.55'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.729'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.377'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8996'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.226'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.676'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1139'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05704'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4628'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.341'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.91%  '
